$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old employee table completely (A1:E6)
$ws.Range("A1:E6").Clear()

# Seed shared strings in the same order as the target workbook:
# Milk(0), Qnty(1), Product(2), Sno(3), TransDate(4)
$ws.Range("B2").Value = "Milk"
$ws.Range("C1").Value = "Qnty"
$ws.Range("B1").Value = "Product"
$ws.Range("A1").Value = "Sno"
$ws.Range("D1").Value = "TransDate"

# Finish row 1 and row 2 numeric/date values
$ws.Range("A2").Value = 1
$ws.Range("C2").Value = 10
$ws.Range("D2").NumberFormat = "mm-dd-yy"
$ws.Range("D2").Value = [DateTime]"2022-01-07"

# Row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Milk"
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = [DateTime]"2022-02-07"

# Make D3's style match D2's exact style entry (avoids creating a duplicate xf)
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Best-fit column D like Excel would for the date column
$ws.Columns("D:D").ColumnWidth = 8.5

$ws.Range("E2").Select()
